$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.518.84'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '1.807.85'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("D5").Value = '228.08'
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").Value = '0.578'
$ws.Range("E6").Value = '  +3.32%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  +6.12%  '
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").Value = '0.0966'
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("D12").Value = '2.068.08'
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("D14").Value = '1.818.24'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Value = '0.652'
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("D16").Value = '4.47'
$ws.Range("E16").Value = '  +3.41%  '
$ws.Range("D17").Value = '34.493.53'
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").Value = '69.96'
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("D19").Value = '245.45'
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("D20").Value = '0.0₃0792'
$ws.Range("E20").Value = '  -1.01%  '
$ws.Range("E21").Value = '  +1.23%  '
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = '4.20'
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("E24").Value = '  +5.03%  '
$ws.Range("D25").Value = '172.44'
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").Value = '8.01'
$ws.Range("E26").Value = '  +8.51%  '
$ws.Range("D27").Value = '16.94'
$ws.Range("E27").Value = '  +1.57%  '
$ws.Range("E28").Value = '  +1.56%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("D32").Value = '0.0530'
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("D34").Value = '1.82'
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("D35").Value = '1.397.18'
$ws.Range("E35").Value = '  -1.49%  '
$ws.Range("D36").Value = '0.677'
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("E37").Value = '  -5.79%  '
$ws.Range("D38").Value = '1.07'
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("E39").Value = '  -0.64%  '
$ws.Range("D40").Value = '83.16'
$ws.Range("E40").Value = '  -3.07%  '
$ws.Range("D41").Value = '0.965'
$ws.Range("E41").Value = '  +1.50%  '
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("D43").Value = '2.43'
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("E44").Value = '  +7.97%  '
$ws.Range("D45").Value = '13.53'
$ws.Range("E45").Value = '  -2.53%  '
$ws.Range("D46").Value = '6.05'
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("D47").Value = '0.0501'
$ws.Range("E47").Value = '  -5.00%  '
$ws.Range("D48").Value = '1.968.84'
$ws.Range("E48").Value = '  -0.75%  '
$ws.Range("D49").Value = '104.55'
$ws.Range("E49").Value = '  -1.14%  '
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("E51").Value = '  -2.69%  '
